# Add 2022-Q3 data
#
# 1) "总计" (summary) sheet: insert a new top data row for 2022-Q3 and push
#    the existing quarters down a row.
# 2) Insert a brand-new worksheet named "2022-Q3" right after "总计" (i.e.
#    before the current "2022-Q2" sheet) and fill it with the per-fund
#    detail row for the new quarter, copying the header/formatting from the
#    "2022-Q2" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "总计" summary sheet - shift rows 2-4 down to 3-5 (bottom-up so we
#    never overwrite data before it has been copied), then write the new
#    2022-Q3 figures into row 2. Range.Copy preserves the existing cell
#    style (s="2" on column A) so we don't need to touch formatting.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

$summary.Range("A4:D4").Copy($summary.Range("A5:D5"))
$summary.Range("A3:D3").Copy($summary.Range("A4:D4"))
$summary.Range("A2:D2").Copy($summary.Range("A3:D3"))

# Fix up the sequential index column that Copy duplicated verbatim.
$summary.Range("A3").Value = 1
$summary.Range("A4").Value = 2
$summary.Range("A5").Value = 3

$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 1
$summary.Range("D2").Value = 0.17

# ---------------------------------------------------------------------
# 2) Brand-new "2022-Q3" worksheet, inserted before "2022-Q2" so the tab
#    order becomes 总计, 2022-Q3, 2022-Q2, 2022-Q1, 2021-Q3.
# ---------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add($wb.Worksheets.Item("2022-Q2"))
$newSheet.Name = "2022-Q3"

$src = $wb.Worksheets.Item("2022-Q2")
$dst = $wb.Worksheets.Item("2022-Q3")

# Copy header row + the row-2 skeleton (values + styles) from 2022-Q2, then
# overwrite with 2022-Q3's own figures.
$src.Range("B1:H1").Copy($dst.Range("B1:H1"))
$src.Range("A2:H2").Copy($dst.Range("A2:H2"))

# Text-like columns must stay text (matches the source data which stores
# these as strings, not numbers) - a leading apostrophe forces Excel to
# keep the numeric-looking input as text, then resetting the style back to
# Normal clears the "quote prefix" flag Excel records for that so the cell
# ends up with no special formatting, same as its neighbours.
$dst.Range("B2").Value = "'398061"
$dst.Range("B2").Style = "Normal"
$dst.Range("C2").Value = "中海消费混合"
$dst.Range("D2").Value = "'3.91"
$dst.Range("D2").Style = "Normal"
$dst.Range("E2").Value = "'85.30"
$dst.Range("E2").Style = "Normal"
$dst.Range("F2").Value = "'4.42"
$dst.Range("F2").Style = "Normal"
$dst.Range("G2").Value = "'0.1728"
$dst.Range("G2").Style = "Normal"
$dst.Range("H2").Value = 7

# Restore the originally active tab (2021-Q3) so the new sheet doesn't
# steal the "selected" flag.
$wb.Worksheets.Item("2021-Q3").Activate()
